$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the error message and timestamp reported on the "report" sheet
# (row 4: Task Code=T1, Name=IPC, description=Monthly CPI, Is Done=No, error=..., date=...)
$ws.Range("F4").Value = "Database failed to get nacional CPI last update date"
$ws.Range("G4").Value = "2022-09-08 22:13:46"
